$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 and row 24 swap values in columns A, B, D, E, F, G, H, Q, R.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

# Read all old values first (using Value2 to get plain scalars, not the
# reflection-describing object that plain .Value seems to yield in this
# runtime) before writing anything, so the swap doesn't clobber itself.
$row23vals = @{}
$row24vals = @{}
foreach ($col in $cols) {
    $row23vals[$col] = $ws.Range("$col" + "23").Value2
    $row24vals[$col] = $ws.Range("$col" + "24").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col" + "23").Value = $row24vals[$col]
    $ws.Range("$col" + "24").Value = $row23vals[$col]
}
